$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.030688913991657
$ws.Cells.Item(2, 4).Value = 1.033777895261852
$ws.Cells.Item(2, 5).Value = 1.03906820358443
$ws.Cells.Item(2, 6).Value = 1.046966858797792
$ws.Cells.Item(2, 9).Value = 1.031885535942883
$ws.Cells.Item(2, 10).Value = 1.035828773694291
$ws.Cells.Item(2, 11).Value = 1.036579266983557
$ws.Cells.Item(2, 12).Value = 1.041854443255225
$ws.Cells.Item(2, 13).Value = 1.049730814449257

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.031741231226777
$ws.Cells.Item(3, 4).Value = 1.034563260479992
$ws.Cells.Item(3, 5).Value = 1.040038460422521
$ws.Cells.Item(3, 6).Value = 1.048106192718484
$ws.Cells.Item(3, 9).Value = 1.032080558505502
$ws.Cells.Item(3, 10).Value = 1.0365221261378
$ws.Cells.Item(3, 11).Value = 1.037173746211646
$ws.Cells.Item(3, 12).Value = 1.042634420362074
$ws.Cells.Item(3, 13).Value = 1.050681047445338

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.032421978830613
$ws.Cells.Item(4, 4).Value = 1.035070716045394
$ws.Cells.Item(4, 5).Value = 1.040666572577727
$ws.Cells.Item(4, 6).Value = 1.04884400588789
$ws.Cells.Item(4, 9).Value = 1.032204654150423
$ws.Cells.Item(4, 10).Value = 1.036970043100419
$ws.Cells.Item(4, 11).Value = 1.03755703308098
$ws.Cells.Item(4, 12).Value = 1.043138792993387
$ws.Cells.Item(4, 13).Value = 1.051295917536339

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.032708124154603
$ws.Cells.Item(5, 4).Value = 1.035283875325801
$ws.Cells.Item(5, 5).Value = 1.04093070032573
$ws.Cells.Item(5, 6).Value = 1.049154323265423
$ws.Cells.Item(5, 9).Value = 1.032256321492679
$ws.Cells.Item(5, 10).Value = 1.037158172758816
$ws.Cells.Item(5, 11).Value = 1.037717835779567
$ws.Cells.Item(5, 12).Value = 1.043350753296416
$ws.Cells.Item(5, 13).Value = 1.051554410195556

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.032756166817261
$ws.Cells.Item(6, 4).Value = 1.035319655422098
$ws.Cells.Item(6, 5).Value = 1.040975052597772
$ws.Cells.Item(6, 6).Value = 1.049206435169524
$ws.Cells.Item(6, 9).Value = 1.032264967188873
$ws.Cells.Item(6, 10).Value = 1.037189750308487
$ws.Cells.Item(6, 11).Value = 1.037744815829531
$ws.Cells.Item(6, 12).Value = 1.04338633778654
$ws.Cells.Item(6, 13).Value = 1.05159781232165

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.032425802479163
$ws.Cells.Item(7, 4).Value = 1.035073564978819
$ws.Cells.Item(7, 5).Value = 1.040670101593319
$ws.Cells.Item(7, 6).Value = 1.048848151810992
$ws.Cells.Item(7, 9).Value = 1.032205346506906
$ws.Cells.Item(7, 10).Value = 1.036972557584637
$ws.Cells.Item(7, 11).Value = 1.037559183036168
$ws.Cells.Item(7, 12).Value = 1.04314162552391
$ws.Cells.Item(7, 13).Value = 1.051299371523062

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.031044585274103
$ws.Cells.Item(8, 4).Value = 1.034043463488221
$ws.Cells.Item(8, 5).Value = 1.039396045977316
$ws.Cells.Item(8, 6).Value = 1.047351780276541
$ws.Cells.Item(8, 9).Value = 1.031951878664176
$ws.Cells.Item(8, 10).Value = 1.036063246080658
$ws.Cells.Item(8, 11).Value = 1.036780459657748
$ws.Cells.Item(8, 12).Value = 1.042118107277319
$ws.Cells.Item(8, 13).Value = 1.050051948869959

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.028609376345163
$ws.Cells.Item(9, 4).Value = 1.032222745322431
$ws.Cells.Item(9, 5).Value = 1.037153238831544
$ws.Cells.Item(9, 6).Value = 1.044719485927107
$ws.Cells.Item(9, 9).Value = 1.031489194396967
$ws.Cells.Item(9, 10).Value = 1.03445535414898
$ws.Cells.Item(9, 11).Value = 1.035397694426084
$ws.Cells.Item(9, 12).Value = 1.040312059179624
$ws.Cells.Item(9, 13).Value = 1.047853879540023

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.026984992996846
$ws.Cells.Item(10, 4).Value = 1.031005240476399
$ws.Cells.Item(10, 5).Value = 1.035659554479717
$ws.Cells.Item(10, 6).Value = 1.042967647965572
$ws.Cells.Item(10, 9).Value = 1.031169976694749
$ws.Cells.Item(10, 10).Value = 1.033379692733725
$ws.Cells.Item(10, 11).Value = 1.034468780174822
$ws.Cells.Item(10, 12).Value = 1.039106369866966
$ws.Cells.Item(10, 13).Value = 1.046388535283523

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.026281395930686
$ws.Cells.Item(11, 4).Value = 1.030477178504828
$ws.Cells.Item(11, 5).Value = 1.03501313509336
$ws.Cells.Item(11, 6).Value = 1.042209798444811
$ws.Cells.Item(11, 9).Value = 1.031029203932895
$ws.Cells.Item(11, 10).Value = 1.032913035308695
$ws.Cells.Item(11, 11).Value = 1.034064877559209
$ws.Cells.Item(11, 12).Value = 1.038583900980823
$ws.Cells.Item(11, 13).Value = 1.045754032877813

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.026020013776208
$ws.Cells.Item(12, 4).Value = 1.030280901869584
$ws.Cells.Item(12, 5).Value = 1.034773079641205
$ws.Cells.Item(12, 6).Value = 1.041928405461056
$ws.Cells.Item(12, 9).Value = 1.030976531931428
$ws.Cells.Item(12, 10).Value = 1.032739564413446
$ws.Cells.Item(12, 11).Value = 1.033914598514627
$ws.Cells.Item(12, 12).Value = 1.038389773031814
$ws.Cells.Item(12, 13).Value = 1.045518350466792

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.026076082723042
$ws.Cells.Item(13, 4).Value = 1.030323009790796
$ws.Cells.Item(13, 5).Value = 1.034824569929548
$ws.Cells.Item(13, 6).Value = 1.041988760415856
$ws.Cells.Item(13, 9).Value = 1.030987847569339
$ws.Cells.Item(13, 10).Value = 1.032776780552771
$ws.Cells.Item(13, 11).Value = 1.033946845262465
$ws.Cells.Item(13, 12).Value = 1.038431416839815
$ws.Cells.Item(13, 13).Value = 1.045568905145192

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.026259790715032
$ws.Cells.Item(14, 4).Value = 1.030460956884818
$ws.Cells.Item(14, 5).Value = 1.034993290943439
$ws.Cells.Item(14, 6).Value = 1.042186536250796
$ws.Cells.Item(14, 9).Value = 1.031024857854761
$ws.Cells.Item(14, 10).Value = 1.03289869887476
$ws.Cells.Item(14, 11).Value = 1.034052460578181
$ws.Cells.Item(14, 12).Value = 1.038567855522897
$ws.Cells.Item(14, 13).Value = 1.045734551296486

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.026372974633472
$ws.Cells.Item(15, 4).Value = 1.030545933309158
$ws.Cells.Item(15, 5).Value = 1.03509725261747
$ws.Cells.Item(15, 6).Value = 1.042308406517207
$ws.Cells.Item(15, 9).Value = 1.03104761040615
$ws.Cells.Item(15, 10).Value = 1.032973799080085
$ws.Cells.Item(15, 11).Value = 1.0341175003243
$ws.Cells.Item(15, 12).Value = 1.038651911981281
$ws.Cells.Item(15, 13).Value = 1.045836611356729

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.02703168355727
$ws.Cells.Item(16, 4).Value = 1.031040267827895
$ws.Cells.Item(16, 5).Value = 1.035702462704827
$ws.Cells.Item(16, 6).Value = 1.0430179588328
$ws.Cells.Item(16, 9).Value = 1.031179265631301
$ws.Cells.Item(16, 10).Value = 1.033410644521688
$ws.Cells.Item(16, 11).Value = 1.034495550559537
$ws.Cells.Item(16, 12).Value = 1.039141036029977
$ws.Cells.Item(16, 13).Value = 1.046430645132447

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.027444812725436
$ws.Cells.Item(17, 4).Value = 1.031350117049331
$ws.Cells.Item(17, 5).Value = 1.036082190286906
$ws.Cells.Item(17, 6).Value = 1.043463231577931
$ws.Cells.Item(17, 9).Value = 1.031261167136932
$ws.Cells.Item(17, 10).Value = 1.033684428137308
$ws.Cells.Item(17, 11).Value = 1.034732242795875
$ws.Cells.Item(17, 12).Value = 1.039447744194719
$ws.Cells.Item(17, 13).Value = 1.04680326706612

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.02768576198494
$ws.Cells.Item(18, 4).Value = 1.031530762467252
$ws.Cells.Item(18, 5).Value = 1.036303713288071
$ws.Cells.Item(18, 6).Value = 1.043723019955628
$ws.Cells.Item(18, 9).Value = 1.031308692901478
$ws.Cells.Item(18, 10).Value = 1.033844035719236
$ws.Cells.Item(18, 11).Value = 1.034870139586725
$ws.Cells.Item(18, 12).Value = 1.039626603531484
$ws.Cells.Item(18, 13).Value = 1.047020611275255

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.027767915742266
$ws.Cells.Item(19, 4).Value = 1.031592343558395
$ws.Cells.Item(19, 5).Value = 1.036379252728077
$ws.Cells.Item(19, 6).Value = 1.043811612675361
$ws.Cells.Item(19, 9).Value = 1.031324856224625
$ws.Cells.Item(19, 10).Value = 1.033898443231947
$ws.Cells.Item(19, 11).Value = 1.034917131378947
$ws.Cells.Item(19, 12).Value = 1.039687583439695
$ws.Cells.Item(19, 13).Value = 1.047094720062507

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.027400490154437
$ws.Cells.Item(20, 4).Value = 1.031316881885773
$ws.Cells.Item(20, 5).Value = 1.03604144556005
$ws.Cells.Item(20, 6).Value = 1.0434154509617
$ws.Cells.Item(20, 9).Value = 1.031252405323061
$ws.Cells.Item(20, 10).Value = 1.033655062617753
$ws.Cells.Item(20, 11).Value = 1.03470686469121
$ws.Cells.Item(20, 12).Value = 1.039414841277743
$ws.Cells.Item(20, 13).Value = 1.04676328824571

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.026205694257544
$ws.Cells.Item(21, 4).Value = 1.030420338513968
$ws.Cells.Item(21, 5).Value = 1.034943605333275
$ws.Cells.Item(21, 6).Value = 1.042128293271012
$ws.Cells.Item(21, 9).Value = 1.031013969811108
$ws.Cells.Item(21, 10).Value = 1.032862800654307
$ws.Cells.Item(21, 11).Value = 1.034021366444235
$ws.Cells.Item(21, 12).Value = 1.038527679380795
$ws.Cells.Item(21, 13).Value = 1.04568577263398

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.025454276586431
$ws.Cells.Item(22, 4).Value = 1.029855888940357
$ws.Cells.Item(22, 5).Value = 1.034253659098087
$ws.Cells.Item(22, 6).Value = 1.04131961969919
$ws.Cells.Item(22, 9).Value = 1.030861841721372
$ws.Cells.Item(22, 10).Value = 1.032363901432299
$ws.Cells.Item(22, 11).Value = 1.033588910762011
$ws.Cells.Item(22, 12).Value = 1.03796953981927
$ws.Cells.Item(22, 13).Value = 1.045008295911444

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.025852636656291
$ws.Cells.Item(23, 4).Value = 1.030155185962224
$ws.Cells.Item(23, 5).Value = 1.034619383305418
$ws.Cells.Item(23, 6).Value = 1.04174825479679
$ws.Cells.Item(23, 9).Value = 1.030942697488677
$ws.Cells.Item(23, 10).Value = 1.032628450591578
$ws.Cells.Item(23, 11).Value = 1.03381830160247
$ws.Cells.Item(23, 12).Value = 1.038265452937314
$ws.Cells.Item(23, 13).Value = 1.045367439186827

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.027420517684113
$ws.Cells.Item(24, 4).Value = 1.031331899686263
$ws.Cells.Item(24, 5).Value = 1.036059856242069
$ws.Cells.Item(24, 6).Value = 1.04343704075417
$ws.Cells.Item(24, 9).Value = 1.031256365169647
$ws.Cells.Item(24, 10).Value = 1.033668331896887
$ws.Cells.Item(24, 11).Value = 1.034718332464108
$ws.Cells.Item(24, 12).Value = 1.039429708807974
$ws.Cells.Item(24, 13).Value = 1.046781352953167

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.029239095401842
$ws.Cells.Item(25, 4).Value = 1.032694098159833
$ws.Cells.Item(25, 5).Value = 1.037732791200228
$ws.Cells.Item(25, 6).Value = 1.045399463707819
$ws.Cells.Item(25, 9).Value = 1.03161070788088
$ws.Cells.Item(25, 10).Value = 1.034871691558839
$ws.Cells.Item(25, 11).Value = 1.035756420508819
$ws.Cells.Item(25, 12).Value = 1.04077925863368
$ws.Cells.Item(25, 13).Value = 1.048422127087495
